$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 2-8 with new forecast-error values
$ws.Range("B2").Value = 0.08631603385188587
$ws.Range("C2").Value = 0.2336348552713851
$ws.Range("D2").Value = 0.06595001910092338
$ws.Range("E2").Value = 0.2568073579571337
$ws.Range("F2").Value = 0.2509970670398742
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.1074344412111381
$ws.Range("C3").Value = 0.2042603668832133
$ws.Range("D3").Value = 0.06469647195393771
$ws.Range("E3").Value = 0.2543550116548477
$ws.Range("F3").Value = 0.2399663286002585
$ws.Range("G3").Value = 13

$ws.Range("B4").Value = 0.1080779521229239
$ws.Range("C4").Value = 0.3459062368002563
$ws.Range("D4").Value = 0.1618177455352025
$ws.Range("E4").Value = 0.4022657648063063
$ws.Range("F4").Value = 0.4047044737270317
$ws.Range("G4").Value = 12

$ws.Range("B5").Value = 0.1348047727392753
$ws.Range("C5").Value = 0.2404736632640445
$ws.Range("D5").Value = 0.1336825235075403
$ws.Range("E5").Value = 0.3656262073587455
$ws.Range("F5").Value = 0.3564564719985849
$ws.Range("G5").Value = 11

$ws.Range("B6").Value = 0.1094716894729319
$ws.Range("C6").Value = 0.3502298540383924
$ws.Range("D6").Value = 0.1684591240845271
$ws.Range("E6").Value = 0.4104377225408589
$ws.Range("F6").Value = 0.416966656392023
$ws.Range("G6").Value = 10

$ws.Range("B7").Value = 0.05568622347126046
$ws.Range("C7").Value = 0.2806577130811624
$ws.Range("D7").Value = 0.1125884625785131
$ws.Range("E7").Value = 0.3355420429372646
$ws.Range("F7").Value = 0.3509607463531708
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = -0.06433122593588818
$ws.Range("C8").Value = 0.3695289678539498
$ws.Range("D8").Value = 0.164134063030832
$ws.Range("E8").Value = 0.4051346233424539
$ws.Range("F8").Value = 0.4381719613125665
$ws.Range("G8").Value = 6

# Row 9 (Q7): now also has an F value, and updated values
$ws.Range("B9").Value = -0.06392244592618833
$ws.Range("C9").Value = 0.1995501394904485
$ws.Range("D9").Value = 0.05897833981413506
$ws.Range("E9").Value = 0.2428545651498754
$ws.Range("F9").Value = 0.2869466694029099
$ws.Range("G9").Value = 3

# New row 10 (Q8)
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -0.06520887812495521
$ws.Range("C10").Value = 0.06520887812495521
$ws.Range("D10").Value = 0.004252197786315262
$ws.Range("E10").Value = 0.06520887812495521
$ws.Range("G10").Value = 1
